# Natmi following Dr Hou advice:
# Recomputed the Efna4-Epha3 ligand/receptor edge statistics (ECs/FAPs/sCs
# senders now cross every target cluster, including the new "M2" cluster)
# and rewrote the corresponding rows on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna4"
$ws.Cells.Item(2, 3).Value = "Epha3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.828516
$ws.Cells.Item(2, 8).Value = 2.485548
$ws.Cells.Item(2, 9).Value = 0.4625620436231038
$ws.Cells.Item(2, 10).Value = 0.4821955800271095
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.030493
$ws.Cells.Item(2, 14).Value = 0.091479
$ws.Cells.Item(2, 15).Value = 0.001018637778668347
$ws.Cells.Item(2, 16).Value = 0.001021037349570579
$ws.Cells.Item(2, 17).Value = 0.025263938388
$ws.Cells.Item(2, 18).Value = 0.227375445492
$ws.Cells.Item(2, 19).Value = 0.0004711831726125295
$ws.Cells.Item(2, 20).Value = 0.0004923396970055281

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna4"
$ws.Cells.Item(3, 3).Value = "Epha3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.828516
$ws.Cells.Item(3, 8).Value = 2.485548
$ws.Cells.Item(3, 9).Value = 0.4625620436231038
$ws.Cells.Item(3, 10).Value = 0.4821955800271095
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 29.65321533333333
$ws.Cells.Item(3, 14).Value = 88.95964599999999
$ws.Cells.Item(3, 15).Value = 0.990584245483253
$ws.Cells.Item(3, 16).Value = 0.9929177316168408
$ws.Cells.Item(3, 17).Value = 24.568163355112
$ws.Cells.Item(3, 18).Value = 221.113470196008
$ws.Cells.Item(3, 19).Value = 0.4582066729715838
$ws.Cells.Item(3, 20).Value = 0.4787805415161844

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna4"
$ws.Cells.Item(4, 3).Value = "Epha3"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.828516
$ws.Cells.Item(4, 8).Value = 2.485548
$ws.Cells.Item(4, 9).Value = 0.4625620436231038
$ws.Cells.Item(4, 10).Value = 0.4821955800271095
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.04031433333333333
$ws.Cells.Item(4, 14).Value = 0.120943
$ws.Cells.Item(4, 15).Value = 0.001346725574891351
$ws.Cells.Item(4, 16).Value = 0.001349898011227873
$ws.Cells.Item(4, 17).Value = 0.033401070196
$ws.Cells.Item(4, 18).Value = 0.300609631764
$ws.Cells.Item(4, 19).Value = 0.0006229441341212426
$ws.Cells.Item(4, 20).Value = 0.0006509148545014657

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Efna4"
$ws.Cells.Item(5, 3).Value = "Epha3"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.828516
$ws.Cells.Item(5, 8).Value = 2.485548
$ws.Cells.Item(5, 9).Value = 0.4625620436231038
$ws.Cells.Item(5, 10).Value = 0.4821955800271095
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.211054
$ws.Cells.Item(5, 14).Value = 0.422108
$ws.Cells.Item(5, 15).Value = 0.007050391163187267
$ws.Cells.Item(5, 16).Value = 0.00471133302236074
$ws.Cells.Item(5, 17).Value = 0.174861615864
$ws.Cells.Item(5, 18).Value = 1.049169695184
$ws.Cells.Item(5, 19).Value = 0.003261243344786174
$ws.Cells.Item(5, 20).Value = 0.002271783959418112

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efna4"
$ws.Cells.Item(6, 3).Value = "Epha3"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.7438396666666667
$ws.Cells.Item(6, 8).Value = 2.231519
$ws.Cells.Item(6, 9).Value = 0.4152870872032183
$ws.Cells.Item(6, 10).Value = 0.4329140288365043
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.030493
$ws.Cells.Item(6, 14).Value = 0.091479
$ws.Cells.Item(6, 15).Value = 0.001018637778668347
$ws.Cells.Item(6, 16).Value = 0.001021037349570579
$ws.Cells.Item(6, 17).Value = 0.02268190295566667
$ws.Cells.Item(6, 18).Value = 0.204137126601
$ws.Cells.Item(6, 19).Value = 0.0004230271160183345
$ws.Cells.Item(6, 20).Value = 0.0004420213925951456

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna4"
$ws.Cells.Item(7, 3).Value = "Epha3"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.7438396666666667
$ws.Cells.Item(7, 8).Value = 2.231519
$ws.Cells.Item(7, 9).Value = 0.4152870872032183
$ws.Cells.Item(7, 10).Value = 0.4329140288365043
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 29.65321533333333
$ws.Cells.Item(7, 14).Value = 88.95964599999999
$ws.Cells.Item(7, 15).Value = 0.990584245483253
$ws.Cells.Item(7, 16).Value = 0.9929177316168408
$ws.Cells.Item(7, 17).Value = 22.05723780914155
$ws.Cells.Item(7, 18).Value = 198.515140282274
$ws.Cells.Item(7, 19).Value = 0.4113768459361379
$ws.Cells.Item(7, 20).Value = 0.4298480154974494

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Efna4"
$ws.Cells.Item(8, 3).Value = "Epha3"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.7438396666666667
$ws.Cells.Item(8, 8).Value = 2.231519
$ws.Cells.Item(8, 9).Value = 0.4152870872032183
$ws.Cells.Item(8, 10).Value = 0.4329140288365043
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.04031433333333333
$ws.Cells.Item(8, 14).Value = 0.120943
$ws.Cells.Item(8, 15).Value = 0.001346725574891351
$ws.Cells.Item(8, 16).Value = 0.001349898011227873
$ws.Cells.Item(8, 17).Value = 0.02998740026855556
$ws.Cells.Item(8, 18).Value = 0.269886602417
$ws.Cells.Item(8, 19).Value = 0.0005592777412587088
$ws.Cells.Item(8, 20).Value = 0.0005843897865590431

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Efna4"
$ws.Cells.Item(9, 3).Value = "Epha3"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.7438396666666667
$ws.Cells.Item(9, 8).Value = 2.231519
$ws.Cells.Item(9, 9).Value = 0.4152870872032183
$ws.Cells.Item(9, 10).Value = 0.4329140288365043
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.211054
$ws.Cells.Item(9, 14).Value = 0.422108
$ws.Cells.Item(9, 15).Value = 0.007050391163187267
$ws.Cells.Item(9, 16).Value = 0.00471133302236074
$ws.Cells.Item(9, 17).Value = 0.1569903370086667
$ws.Cells.Item(9, 18).Value = 0.9419420220520001
$ws.Cells.Item(9, 19).Value = 0.00292793640980335
$ws.Cells.Item(9, 20).Value = 0.002039602159900652

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Efna4"
$ws.Cells.Item(10, 3).Value = "Epha3"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.21879
$ws.Cells.Item(10, 8).Value = 0.43758
$ws.Cells.Item(10, 9).Value = 0.1221508691736778
$ws.Cells.Item(10, 10).Value = 0.08489039113638626
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.030493
$ws.Cells.Item(10, 14).Value = 0.091479
$ws.Cells.Item(10, 15).Value = 0.001018637778668347
$ws.Cells.Item(10, 16).Value = 0.001021037349570579
$ws.Cells.Item(10, 17).Value = 0.00667156347
$ws.Cells.Item(10, 18).Value = 0.04002938082
$ws.Cells.Item(10, 19).Value = 0.0001244274900374831
$ws.Cells.Item(10, 20).Value = 0.000086676259969905623554292662

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Efna4"
$ws.Cells.Item(11, 3).Value = "Epha3"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.21879
$ws.Cells.Item(11, 8).Value = 0.43758
$ws.Cells.Item(11, 9).Value = 0.1221508691736778
$ws.Cells.Item(11, 10).Value = 0.08489039113638626
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 29.65321533333333
$ws.Cells.Item(11, 14).Value = 88.95964599999999
$ws.Cells.Item(11, 15).Value = 0.990584245483253
$ws.Cells.Item(11, 16).Value = 0.9929177316168408
$ws.Cells.Item(11, 17).Value = 6.487826982779999
$ws.Cells.Item(11, 18).Value = 38.92696189667999
$ws.Cells.Item(11, 19).Value = 0.1210007265755312
$ws.Cells.Item(11, 20).Value = 0.08428917460320702

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Efna4"
$ws.Cells.Item(12, 3).Value = "Epha3"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.21879
$ws.Cells.Item(12, 8).Value = 0.43758
$ws.Cells.Item(12, 9).Value = 0.1221508691736778
$ws.Cells.Item(12, 10).Value = 0.08489039113638626
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.04031433333333333
$ws.Cells.Item(12, 14).Value = 0.120943
$ws.Cells.Item(12, 15).Value = 0.001346725574891351
$ws.Cells.Item(12, 16).Value = 0.001349898011227873
$ws.Cells.Item(12, 17).Value = 0.00882037299
$ws.Cells.Item(12, 18).Value = 0.05292223793999999
$ws.Cells.Item(12, 19).Value = 0.0001645036995113995
$ws.Cells.Item(12, 20).Value = 0.0001145933701673641

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Efna4"
$ws.Cells.Item(13, 3).Value = "Epha3"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.21879
$ws.Cells.Item(13, 8).Value = 0.43758
$ws.Cells.Item(13, 9).Value = 0.1221508691736778
$ws.Cells.Item(13, 10).Value = 0.08489039113638626
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.211054
$ws.Cells.Item(13, 14).Value = 0.422108
$ws.Cells.Item(13, 15).Value = 0.007050391163187267
$ws.Cells.Item(13, 16).Value = 0.00471133302236074
$ws.Cells.Item(13, 17).Value = 0.04617650466
$ws.Cells.Item(13, 18).Value = 0.18470601864
$ws.Cells.Item(13, 19).Value = 0.0008612114085977422
$ws.Cells.Item(13, 20).Value = 0.000399946903041976
